$wb = $excel.ActiveWorkbook

# --- Sheet1: MAIN_CONTROLLER ---
$ws1 = $wb.Worksheets.Item("MAIN_CONTROLLER")
$ws1.Range("B2").Value = "N"
$ws1.Range("B6").Value = "Y"
$ws1.Range("D6").Value = "FOS"
$ws1.Range("E6").Value = "FOS"
$ws1.Range("B6").Select() | Out-Null

# --- Sheet2: DATASHEET ---
$ws2 = $wb.Worksheets.Item("DATASHEET")
# NOTE: update D5 before D3 so the shared-string table ends up built in
# the same order as it would be if the two cells were edited that order.
$ws2.Range("D5").Value = "FOS3UW_to_postSanction2.xlsx"
$ws2.Range("D3").Value = "FOS7_AfterPostSanction.xlsx"

# New row 7 - duplicate of row 6 (CPC_1stTouchPoint_Approval), carry over formatting
$ws2.Range("A6:F6").Copy() | Out-Null
$ws2.Range("A7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws2.Range("A7").Value = 2
$ws2.Range("B7").Value = "Y"
$ws2.Range("C7").Value = "CPC_1stTouchPoint_Approval"
$ws2.Range("D7").Value = "CPC_Dynamic.xlsx"
$ws2.Range("E7").Value = 10
$ws2.Range("F7").Value = 20

# New row 8 - FOS / FOS7_AfterPostSanction.xlsx, no special formatting
$ws2.Range("A8").Value = 2
$ws2.Range("B8").Value = "Y"
$ws2.Range("C8").Value = "FOS"
$ws2.Range("D8").Value = "FOS7_AfterPostSanction.xlsx"
$ws2.Range("E8").Value = 10
$ws2.Range("F8").Value = 20

$ws2.Range("D15").Select() | Out-Null

# --- Sheet3: MOBILE_CONFIGURATION ---
$ws3 = $wb.Worksheets.Item("MOBILE_CONFIGURATION")
# Copy formatting from row4 into the new row5 before filling values
$ws3.Range("A4:P4").Copy() | Out-Null
$ws3.Range("A5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws3.Range("A5").Value = 1
$ws3.Range("B5").Value = "Y"
$ws3.Range("C5").Value = "FOS"
$ws3.Range("D5").Value = "com.mahindra.fospreprod"
$ws3.Range("E5").Value = "com.mahindra.fosbeta.MainActivity"
$ws3.Range("F5").Value = "d4a4d1d2"
$ws3.Range("G5").Value = "Android"
$ws3.Range("H5").Value = 11
$ws3.Range("I5").Value = "No"
$ws3.Range("J5").Value = "Yes"
$ws3.Range("K5").Value = "C:\\Users\\biswa\\Documents\\base.apk"
$ws3.Range("L5").Value = "0.0.0.0:4723"
$ws3.Range("M5").Value = "BrowserStack"
$ws3.Range("N5").Value = "biswajitsahoo_0n9ypv"
$ws3.Range("O5").Value = "qZHZfSFttvdThCVVX6Ki"

$ws3.Range("C5").Select() | Out-Null

$ws1.Activate() | Out-Null
